# "Acutualizacion plan de desarrollo"
# - Rename "Hoja1" to "Plan de desarrollo"
# - Hide the helper lookup sheet "Hoja2"
# - Move the active selection on the main sheet to G3

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Plan de desarrollo"

$ws2 = $wb.Worksheets.Item("Hoja2")
$ws2.Visible = $false

$ws1.Activate()
$ws1.Range("G3").Select()
